# Apply the changes described by the diff:
# - Make "Mes actual emp" (2nd sheet) the active tab instead of "Mes actual cel" (1st sheet)
# - Update the selection on "Mes actual emp" to C3
# - Round the remaining-points value on "Mes actual emp"!C2 to two decimals (84.1688.. -> 84.16)

$wb = $excel.ActiveWorkbook

$wsEmp = $wb.Worksheets.Item("Mes actual emp")

# Round the remaining-points value to two decimals.
$wsEmp.Range("C2").Value = 84.16

# Activate "Mes actual emp" so it becomes the selected/active tab for the workbook
# (this also clears tabSelected on the previously active sheet).
$wsEmp.Activate()

# Update the selected cell on the now-active sheet.
$wsEmp.Range("C3").Select()
